$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Music/Name entry in B5 from "Seth Harmon" to "Drew Smith"
$ws.Range("B5").Value = "Drew Smith"

# Move the active selection to D12 (matches saved cursor position)
$ws.Range("D12").Select()
